# GopherLife.pptx edit: bump header/footer date fields to 26/01/2019,
# rename the "Problem 5 - Optimization" slides to "Part 5" / "Part 4",
# trim two paragraphs + resize/reflow slide 23's content placeholder and
# picture.

$p = $ppt.ActivePresentation

function Set-DatePlaceholderText {
    param($shapes)
    for ($shapeIdx = 1; $shapeIdx -le $shapes.Count; $shapeIdx++) {
        $sh = $shapes.Item($shapeIdx)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = "26/01/2019"
        }
    }
}

# 1. Slide master date placeholder.
Set-DatePlaceholderText $p.SlideMaster.Shapes

# 2. Every slide layout's date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($layoutIdx = 1; $layoutIdx -le $layouts.Count; $layoutIdx++) {
    Set-DatePlaceholderText $layouts.Item($layoutIdx).Shapes
}

# 3. Notes master date placeholder.
Set-DatePlaceholderText $p.NotesMaster.Shapes

# 4. Slide 22 title: "Problem 5 - Optimization" -> "Part 5 - Optimization".
$s22 = $p.Slides.Item(22)
$s22.Shapes.Item("Title 1").TextFrame.TextRange.Text = "Part 5 - Optimization"

# 5. Slide 23 title: "Problem 5 - Optimization" -> "Part 4 - Optimization".
$s23 = $p.Slides.Item(23)
$s23.Shapes.Item("Title 1").TextFrame.TextRange.Text = "Part 4 - Optimization"

# 6. Slide 23 content placeholder: widen it and drop the two
#    "(Note this should show an image of pre-spiral results)" /
#    "What was the cause?" paragraphs.
$content = $s23.Shapes.Item("Content Placeholder 2")
$content.Width = 5568412

$cr = [char]13
$newBody = (
    "At only 6% of my target goal it was taking nearly " + [char]0xBD + " a second to process a single frame. " + $cr +
    $cr +
    "Is this just how fast it takes to process things? " + $cr +
    $cr +
    "Was my goal un-reasonable?" + $cr
)
$content.TextFrame.TextRange.Text = $newBody

# 7. Slide 23 picture: reposition (size unchanged).
$pic = $s23.Shapes.Item("Picture 3")
$pic.Left = 5973329
$pic.Top = 2052916
